# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" worksheet (fund-holdings detail) right after
#   "2021-Q4" and before "总计".
# - Update the "总计" (totals) sheet to add a leading 2022-Q1 row (and
#   shift the rest down).
#
# Implementation note: rather than inserting a brand-new blank sheet for
# "2022-Q1" and separately retyping the existing "总计" rows, we repurpose
# the worksheet object that is currently last ("总计") to become
# "2022-Q1" and rewrite its data, then append a *fresh* worksheet named
# "总计" at the end with the refreshed totals table. That mirrors the
# sheetId/rId bookkeeping Excel itself produces for this kind of edit
# (2022-Q1 keeps sheetId 5 / rId5, the recreated 总计 gets sheetId 6 /
# rId6) and matches the target sheet order.

$wb = $excel.ActiveWorkbook

# Helper style source: "2021-Q4" already has the exact header/index-column
# formatting ("基金代码"/"基金名称"/... table) that "2022-Q1" needs.
$q4Sheet = $wb.Worksheets.Item(4)

# The sheet that is currently last ("总计") becomes "2022-Q1".
$newQ1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$newQ1.Cells.Clear()
$newQ1.Name = "2022-Q1"

# Copy the formatted header row (B1:H1) and index column (A2:A3) from
# "2021-Q4" so the new sheet picks up the workbook's existing shared
# style (bold, centered, thin border) instead of minting a new one.
$q4Sheet.Range("B1:H1").Copy()
$newQ1.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2:A3").Copy()
$newQ1.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newQ1.Range("B1").Value = "基金代码"
$newQ1.Range("C1").Value = "基金名称"
$newQ1.Range("D1").Value = "基金规模"
$newQ1.Range("E1").Value = "股票总仓位"
$newQ1.Range("F1").Value = "仓位占比"
$newQ1.Range("G1").Value = "持有市值(亿元)"
$newQ1.Range("H1").Value = "仓位排名"

$newQ1.Range("A2").Value = 0
$newQ1.Range("A3").Value = 1

# Text-typed data cells (match source workbook, which stores these as
# plain strings, not numbers). NumberFormat "@" forces text entry so
# numeric-looking strings ("1.29", "010375", ...) aren't silently
# coerced into numbers; ClearFormats afterwards drops back to the
# sheet's default (unstyled) cell so no stray style index is left
# behind, matching the rest of the table's plain data cells.
$textCells = @(
    @{ Addr = "B2"; Val = "010375" },
    @{ Addr = "C2"; Val = "国金鑫悦经济新动能混合A" },
    @{ Addr = "D2"; Val = "1.29" },
    @{ Addr = "E2"; Val = "94.93" },
    @{ Addr = "F2"; Val = "4.43" },
    @{ Addr = "G2"; Val = "0.0571" },
    @{ Addr = "B3"; Val = "010376" },
    @{ Addr = "C3"; Val = "国金鑫悦经济新动能混合C" },
    @{ Addr = "D3"; Val = "0.28" },
    @{ Addr = "E3"; Val = "94.93" },
    @{ Addr = "F3"; Val = "4.43" },
    @{ Addr = "G3"; Val = "0.0124" }
)
foreach ($cell in $textCells) {
    $rng = $newQ1.Range($cell.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cell.Val
    $rng.ClearFormats()
}

$newQ1.Range("H2").Value = 9
$newQ1.Range("H3").Value = 9

# Recreate "总计" as a brand-new trailing sheet with the refreshed totals
# table (2022-Q1 row added on top, everything else shifted down one).
$total = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$total.Name = "总计"

$newQ1.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$newQ1.Range("A2:A3").Copy()
$total.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.07000000000000001

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 8
$total.Range("D3").Value = 0.4

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.21

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.01

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2020-Q4"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.08
